$wb = $excel.ActiveWorkbook

# Sheet 1: 展览 (Exhibition)
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("F2").Value = 226
$ws1.Range("F5").Value = 23
$ws1.Range("F6").Value = 83
$ws1.Range("F8").Value = 380
$ws1.Range("F9").Value = 4601
$ws1.Range("F10").Value = 4601
$ws1.Range("F12").Value = 445
$ws1.Range("F13").Value = 1080
$ws1.Range("F15").Value = 4092
$ws1.Range("F17").Value = 157
$ws1.Range("F18").Value = 43
$ws1.Range("F19").Value = 199
$ws1.Range("F20").Value = 3388
$ws1.Range("F24").Value = 2935
$ws1.Range("F25").Value = 119
$ws1.Range("F28").Value = 139
$ws1.Range("F29").Value = 170
$ws1.Range("F30").Value = 171
$ws1.Range("F34").Value = 51
$ws1.Range("F36").Value = 5309
$ws1.Range("F37").Value = 725
$ws1.Range("F38").Value = 381
$ws1.Range("F41").Value = 13
$ws1.Range("F42").Value = 1043
$ws1.Range("F43").Value = 427
$ws1.Range("F45").Value = 1928
$ws1.Range("F46").Value = 294
$ws1.Range("F48").Value = 685
$ws1.Range("F49").Value = 820

# Sheet 2: 演出 (Performance)
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("F5").Value = 78
$ws2.Range("F8").Value = 58
$ws2.Range("F19").Value = 37
$ws2.Range("F21").Value = 716

# Sheet 4: 全部类型 (All types)
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("F5").Value = 226
$ws4.Range("F6").Value = 23
$ws4.Range("F7").Value = 78
$ws4.Range("F9").Value = 380
$ws4.Range("F10").Value = 4601
$ws4.Range("F11").Value = 4601
$ws4.Range("F14").Value = 58
$ws4.Range("F16").Value = 445
$ws4.Range("F17").Value = 1080
$ws4.Range("F19").Value = 4092
$ws4.Range("F21").Value = 157
$ws4.Range("F22").Value = 199
$ws4.Range("F23").Value = 3388
$ws4.Range("F24").Value = 2935
$ws4.Range("F25").Value = 119
$ws4.Range("F27").Value = 139
$ws4.Range("F28").Value = 170
$ws4.Range("F29").Value = 171
$ws4.Range("F33").Value = 51
$ws4.Range("F36").Value = 5309
$ws4.Range("F38").Value = 725
$ws4.Range("F39").Value = 381
$ws4.Range("F43").Value = 1043
$ws4.Range("F44").Value = 427
$ws4.Range("F46").Value = 1928
$ws4.Range("F47").Value = 294
$ws4.Range("F49").Value = 685
$ws4.Range("F50").Value = 820
